$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B (Coin name) swaps for rows 7/8 ---
$ws.Range("B7").Value = "USDC"
$ws.Range("B8").Value = "Solana"

# --- Column C (Link) swaps for rows 7/8 ---
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("C8").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"

# --- Column D (Price) updates ---
$ws.Range("D2").Value = "37.940.61"
$ws.Range("D3").Value = "2.114.77"
$ws.Range("D5").Value = "'234.65"
$ws.Range("D7").Value = "'1.00"
$ws.Range("D8").Value = "'57.71"
$ws.Range("D9").Value = "'0.391"
$ws.Range("D12").Value = "2.432.11"
$ws.Range("D13").Value = "'14.46"
$ws.Range("D14").Value = "'21.23"
$ws.Range("D15").Value = "'0.782"
$ws.Range("D17").Value = "2.116.96"
$ws.Range("D18").Value = "37.862.29"
$ws.Range("D19").Value = "'6.22"
$ws.Range("D20").Value = "'70.41"
$ws.Range("D22").Value = "'227.53"
$ws.Range("D24").Value = "'2.41"
$ws.Range("D25").Value = "'2.40"
$ws.Range("D26").Value = "'169.35"
$ws.Range("D27").Value = "'0.139"
$ws.Range("D28").Value = "'8.96"
$ws.Range("D29").Value = "'1.42"
$ws.Range("D30").Value = "'19.56"
$ws.Range("D31").Value = "'0.119"
$ws.Range("D32").Value = "'4.62"
$ws.Range("D34").Value = "'0.0624"
$ws.Range("D35").Value = "'4.55"
$ws.Range("D36").Value = "'3.45"
$ws.Range("D40").Value = "'0.1000"
$ws.Range("D42").Value = "'96.92"
$ws.Range("D43").Value = "1.462.99"
$ws.Range("D47").Value = "'4.10"
$ws.Range("D48").Value = "'15.55"
$ws.Range("D51").Value = "2.316.15"

# --- Column E (Volume 1h) updates ---
$ws.Range("E2").Value = "  +2.07%  "
$ws.Range("E3").Value = "  +3.03%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E5").Value = "  +1.22%  "
$ws.Range("E6").Value = "  +1.20%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +0.66%  "
$ws.Range("E9").Value = "  +2.90%  "
$ws.Range("E10").Value = "  +3.37%  "
$ws.Range("E11").Value = "  +1.55%  "
$ws.Range("E12").Value = "  +3.15%  "
$ws.Range("E13").Value = "  +1.39%  "
$ws.Range("E14").Value = "  +2.54%  "
$ws.Range("E15").Value = "  +1.47%  "
$ws.Range("E16").Value = "  +1.67%  "
$ws.Range("E17").Value = "  +2.92%  "
$ws.Range("E18").Value = "  +2.02%  "
$ws.Range("E19").Value = "  -1.26%  "
$ws.Range("E20").Value = "  +2.16%  "
$ws.Range("E21").Value = "  +2.14%  "
$ws.Range("E22").Value = "  +1.62%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("E25").Value = "  +1.76%  "
$ws.Range("E26").Value = "  +2.39%  "
$ws.Range("E27").Value = "  +10.61%  "
$ws.Range("E28").Value = "  +2.58%  "
$ws.Range("E29").Value = "  -1.35%  "
$ws.Range("E30").Value = "  +3.16%  "
$ws.Range("E31").Value = "  +1.90%  "
$ws.Range("E32").Value = "  +4.30%  "
$ws.Range("E33").Value = "  +3.77%  "
$ws.Range("E34").Value = "  +2.41%  "
$ws.Range("E35").Value = "  +0.95%  "
$ws.Range("E36").Value = "  +6.06%  "
$ws.Range("E37").Value = "  +5.27%  "
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("E39").Value = "  -5.50%  "
$ws.Range("E40").Value = "  +8.04%  "
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("E42").Value = "  +1.25%  "
$ws.Range("E43").Value = "  -1.07%  "
$ws.Range("E44").Value = "  +2.25%  "
$ws.Range("E45").Value = "  -0.64%  "
$ws.Range("E46").Value = "  +4.41%  "
$ws.Range("E47").Value = "  -11.38%  "
$ws.Range("E48").Value = "  +2.39%  "
$ws.Range("E49").Value = "  +3.69%  "
$ws.Range("E50").Value = "  +2.15%  "
$ws.Range("E51").Value = "  +3.12%  "
